# Weathermap forecast update for RAAL - GHI_2024-02-04.xlsx
#
# Updates:
#  - sunrise/sunset timestamps (shared across "Daily" and "Hourly" sheets)
#  - longitude value (column B) on both sheets
#  - Daily sheet clear/cloudy sky daily totals (row 2)
#  - Hourly sheet clear/cloudy sky hourly values (rows 9-19, i.e. hours 7-17)

$wb = $excel.ActiveWorkbook

$daily  = $wb.Worksheets.Item("Daily")
$hourly = $wb.Worksheets.Item("Hourly")

$newSunrise = "2024-02-04T07:40:20"
$newSunset  = "2024-02-04T17:30:04"
$newLon     = 24.724419

# ---- Daily sheet (row 2) ----
$daily.Range("B2").Value = $newLon
$daily.Range("E2").Value = $newSunrise
$daily.Range("F2").Value = $newSunset

$daily.Range("G2").Value = 2618.52
$daily.Range("H2").Value = 5818.4
$daily.Range("I2").Value = 681.5700000000001
$daily.Range("J2").Value = 654.63
$daily.Range("L2").Value = 654.63

# ---- Hourly sheet ----
# Update longitude and sunrise/sunset text for every data row (2-25)
for ($r = 2; $r -le 25; $r++) {
    $hourly.Range("B$r").Value = $newLon
    $hourly.Range("E$r").Value = $newSunrise
    $hourly.Range("F$r").Value = $newSunset
}

# Hour-by-hour clear/cloudy sky values that changed (rows 9-19 => hours 7-17)
$hourlyUpdates = @{
    9  = @{ H = 1.69;                 I = 15.57;               J = 3.29;                 K = 0.42;                 M = 0.42 }
    10 = @{ H = 79.29000000000001;    I = 362.64;               J = 41.94;                K = 19.82;                M = 19.82 }
    11 = @{ H = 215.94;               I = 601.54;               J = 68.01000000000001;    K = 53.99;                M = 53.99 }
    12 = @{ H = 335.24;               I = 709.96;               J = 82.52;                K = 83.81;                M = 83.81 }
    13 = @{ H = 415.56;               I = 762.61;               J = 90.41;                K = 103.89;               M = 103.89 }
    14 = @{ H = 446.56;               I = 780.15;               J = 93.18000000000001;    K = 111.64;               M = 111.64 }
    15 = @{ H = 424.64;               I = 767.9400000000001;    J = 91.22;                K = 106.16;               M = 106.16 }
    16 = @{ H = 352.33;               I = 722.3;                J = 84.27;                K = 88.08;                M = 88.08 }
    17 = @{ H = 238.69;               I = 626.37;               J = 71.11;                K = 59.67;                M = 59.67 }
    18 = @{ H = 102.46;               I = 421.89;               J = 47.82;                K = 25.61;                M = 25.61 }
    19 = @{ H = 6.12;                 I = 47.43;                J = 7.8;                  K = 1.53;                 M = 1.53 }
}

foreach ($row in $hourlyUpdates.Keys) {
    $vals = $hourlyUpdates[$row]
    foreach ($col in $vals.Keys) {
        $hourly.Range("$col$row").Value = $vals[$col]
    }
}
